$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert 4 new rows (one into each of the 4 plot-data groups) that add a
# "PrecLand" (precision landing) trace alongside the existing EKF1/EKF2/DES
# traces, and renumber the "Plot No" (col A) for the groups that shift down.
#
# Using ORIGINAL row numbers, insert a blank row before each of these rows
# (processing top-to-bottom with a running offset keeps this simple):
#   before orig row 8  -> new row 8  (Vx/VN group)   : PL/vX  -> VN_{PrecLand}
#   before orig row 10 -> new row 11 (PN group)      : PL/pX  -> PN_{PrecLand}
#   before orig row 13 -> new row 15 (Vy/VE group)   : PL/vY  -> VE_{PrecLand}
#   before orig row 15 -> new row 18 (PE group)      : PL/pY  -> PE_{PrecLand}
# ---------------------------------------------------------------------------

$ws.Rows(8).Insert()
$ws.Rows(11).Insert()
$ws.Rows(15).Insert()
$ws.Rows(18).Insert()

# Renumber column A ("Plot No") for the rows that used to read 2/3/4 and now
# read 1/2/3 (the PrecLand rows use the same Plot No as the rest of their
# group).
$ws.Range("A8").Value = 1
$ws.Range("A9").Value = 2
$ws.Range("A10").Value = 2
$ws.Range("A11").Value = 2
$ws.Range("A12").Value = 2
$ws.Range("A13").Value = 3
$ws.Range("A14").Value = 3
$ws.Range("A15").Value = 3
$ws.Range("A16").Value = 4
$ws.Range("A17").Value = 4
$ws.Range("A18").Value = 4
$ws.Range("A19").Value = 4

$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "S"
$ws.Range("E8").Value = "Time [ s ]"
$ws.Range("F8").Value = "Vx~[~m/s~]"
$ws.Range("G8").Value = "Vertical"
$ws.Range("K8").Value = "PL/vX"
$ws.Range("O8").Value = -0.01
$ws.Range("R8").Value = "VN_{PrecLand}"
$ws.Range("S8").Value = "m/s"

$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "S"
$ws.Range("E11").Value = "Time [ s ]"
$ws.Range("F11").Value = "PN~[~m~]"
$ws.Range("G11").Value = "Vertical"
$ws.Range("K11").Value = "PL/pX"
$ws.Range("O11").Value = -0.01
$ws.Range("R11").Value = "PN_{PrecLand}"
$ws.Range("S11").Value = "m"

$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = "S"
$ws.Range("E15").Value = "Time [ s ]"
$ws.Range("F15").Value = "Vy~[~m/s~]"
$ws.Range("G15").Value = "Vertical"
$ws.Range("K15").Value = "PL/vY"
$ws.Range("O15").Value = -0.01
$ws.Range("R15").Value = "VE_{PrecLand}"
$ws.Range("S15").Value = "m/s"

$ws.Range("B18").Value = 2
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = "S"
$ws.Range("E18").Value = "Time [ s ]"
$ws.Range("F18").Value = "PE~[~m~]"
$ws.Range("G18").Value = "Vertical"
$ws.Range("K18").Value = "PL/pY"
$ws.Range("O18").Value = -0.01
$ws.Range("R18").Value = "PE_{PrecLand}"
$ws.Range("S18").Value = "m"

# ---------------------------------------------------------------------------
# Column-width tweak: column L got narrowed (was a single wide spacer column)
# and split visually into three narrow spacer columns L:N.
# ---------------------------------------------------------------------------
$ws.Columns(12).ColumnWidth = 4
$ws.Columns(13).ColumnWidth = 3.8333333333333335
$ws.Columns(14).ColumnWidth = 4.333333333333333

# ---------------------------------------------------------------------------
# Leave the active selection on the last data row, matching the saved
# worksheet view.
# ---------------------------------------------------------------------------
$ws.Rows(18).EntireRow.Select() | Out-Null
